# Generate Report for Handoff
#
# The localization-status report is regenerated: the "be1c3c1d..." source
# file moves from "Handed back: in sync with en-US" to "Ready for handoff"
# (with a new handback timestamp + stale-handback error detail), while the
# "fcce5ca1..." source file keeps the "Handed back: in sync with en-US"
# status it previously had. Because rows are (re)generated sorted by
# status/filename, the two data rows on each per-locale sheet swap places,
# and the two hyperlinks on the Overview sheet swap their display text to
# match.

$wb = $excel.ActiveWorkbook

$be1c = "be1c3c1d-6466-4453-a858-400c0a0e22a2.md"
$fcce = "fcce5ca1-8c39-44ff-a101-cd68a19ab1e9.md"

# ---------------------------------------------------------------------
# Overview sheet: row 3 (fcce5ca1...) status/date cells change, and the
# two hyperlink display strings swap (B2 <-> B3) while keeping the same
# r:id (and therefore the same target URL) they already had.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-01 00:52:44"

foreach ($h in $wsOverview.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$B$2') {
        $h.TextToDisplay = "e2e\" + $fcce
    } elseif ($addr -eq '$B$3') {
        $h.TextToDisplay = "e2e\" + $be1c
    }
}

# ---------------------------------------------------------------------
# Per-locale sheets (zh-cn, de-de): the two data rows swap file identity
# (fcce5ca1... moves into row 2 keeping "Handed back: in sync with
# en-US", be1c3c1d... moves into row 3 and becomes "Ready for handoff"
# with a refreshed handback datetime and a new Error Detail message).
# ---------------------------------------------------------------------
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0b6bec0b8ad6814f5f5c74a790bada42b669a884/e2e/be1c3c1d-6466-4453-a858-400c0a0e22a2.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d478341f0c58fa417e3fa781a08b9d74134467ce/e2e/be1c3c1d-6466-4453-a858-400c0a0e22a2.md."

function Update-LocaleSheet {
    param(
        [string]$sheetName,
        [string]$row2HandoffFile,
        [string]$row3HandoffFile,
        [string]$row2HandbackFile,
        [string]$row3HandbackFile,
        [string]$row2HandoffDate,
        [string]$row3HandbackDate
    )

    $ws = $wb.Worksheets.Item($sheetName)

    # Row 2 -> fcce5ca1..., still "Handed back: in sync with en-US"
    $ws.Range("A2").Value = $fcce
    $ws.Range("C2").Value = "Handed back: in sync with en-US"
    $ws.Range("G2").Value = $row2HandoffFile
    $ws.Range("H2").Value = $row2HandoffDate
    $ws.Range("I2").Value = $fcce
    $ws.Range("J2").Value = $row2HandoffFile
    $ws.Range("K2").Value = $row2HandbackFile

    # Row 3 -> be1c3c1d..., now "Ready for handoff"
    $ws.Range("A3").Value = $be1c
    $ws.Range("C3").Value = "Ready for handoff"
    $ws.Range("G3").Value = $row3HandoffFile
    $ws.Range("H3").Value = $row3HandbackDate
    $ws.Range("I3").Value = $be1c
    $ws.Range("J3").Value = $row3HandoffFile
    $ws.Range("K3").Value = $row3HandbackFile
    $ws.Range("P3").Value = $errorDetail

    foreach ($h in $ws.Hyperlinks) {
        $addr = $h.Range.Address()
        if ($addr -eq '$A$2' -or $addr -eq '$I$2') {
            $h.TextToDisplay = $fcce
        } elseif ($addr -eq '$A$3' -or $addr -eq '$I$3') {
            $h.TextToDisplay = $be1c
        }
    }

    # Excel's ColumnWidth (chars) differs from the stored OOXML width
    # (chars + ~5px padding / Maximum Digit Width) by a small constant;
    # 39.14 round-trips to a stored width of exactly 40.
    $ws.Columns.Item(16).ColumnWidth = 39.14
}

Update-LocaleSheet "zh-cn" `
    "fcce5ca1-8c39-44ff-a101-cd68a19ab1e9.32db8f3d46f70a58a0f17b0322d1bf7b6dac33ff.zh-cn.xlf" `
    "be1c3c1d-6466-4453-a858-400c0a0e22a2.d0bef59ea020746cc19c13e92b16ac1d6a474957.zh-cn.xlf" `
    "2016-09-01 00:52:15" `
    "2016-09-01 00:52:15" `
    "2016-09-01 00:51:46" `
    "2016-09-01 00:52:40"

Update-LocaleSheet "de-de" `
    "fcce5ca1-8c39-44ff-a101-cd68a19ab1e9.32db8f3d46f70a58a0f17b0322d1bf7b6dac33ff.de-de.xlf" `
    "be1c3c1d-6466-4453-a858-400c0a0e22a2.d0bef59ea020746cc19c13e92b16ac1d6a474957.de-de.xlf" `
    "2016-09-01 00:52:22" `
    "2016-09-01 00:52:22" `
    "2016-09-01 00:51:52" `
    "2016-09-01 00:52:44"
